# Add two new columns, I ("I0") and J ("IF"), to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells (row 1) ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy formatting (bold font, border, centered/top alignment) from the
# existing "IP" header cell (H1) onto the two new header cells so they
# reuse the same cell style rather than creating a new one.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data values (rows 2-69) for columns I and J ---
$iValues = @(5,7,6,6,7,9,8,4,7,8,6,6,6,6,6,5,6,9,10,7,6,9,7,6,8,5,6,7,8,7,9,7,6,7,7,6,6,5,6,7,5,7,5,10,10,6,7,7,1,7,7,6,5,5,8,5,8,7,6,6,7,6,6,7,8,6,3,3)
$jValues = @(6,7,6,6,8,9,8,4,7,9,7,6,6,6,6,5,7,9,10,7,7,9,8,6,8,5,7,7,8,8,9,7,7,7,7,6,6,6,6,7,6,8,5,10,10,7,7,8,3,7,7,7,6,6,8,6,8,7,6,8,8,7,7,7,8,7,4,3)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
